# Weekly update: a new week of price data (2021-08-06, serial 44414) is
# prepended to the historical series. This pushes every existing record
# down by one week (2 rows: "Primera" + "Segunda"), so the oldest week
# that used to occupy rows 133:134 now ends up in rows 135:136.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new week by inserting two blank rows at the top of
# the data block (row 36 is the first data row below the header + the
# other "Agrícola del Norte S.A. de Arica" blocks already in the sheet).
$ws.Rows("36:37").Insert()

# Row 36: new week, calidad "Primera"
$ws.Cells.Item(36,1).Value  = 1
$ws.Cells.Item(36,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(36,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(36,4).Value  = 44414
$ws.Cells.Item(36,5).Value  = 15
$ws.Cells.Item(36,6).Value  = 100114014
$ws.Cells.Item(36,7).Value  = "Betarraga"
$ws.Cells.Item(36,8).Value  = "Sin especificar"
$ws.Cells.Item(36,9).Value  = "Primera"
$ws.Cells.Item(36,10).Value = 700
$ws.Cells.Item(36,11).Value = 500
$ws.Cells.Item(36,12).Value = 550
$ws.Cells.Item(36,13).Value = 525
$ws.Cells.Item(36,14).Value = "`$/paquete 4 unidades"
$ws.Cells.Item(36,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(36,16).Value = 131
$ws.Cells.Item(36,17).Value = 4
$ws.Cells.Item(36,18).Value = "Hortaliza"

# Row 37: new week, calidad "Segunda"
$ws.Cells.Item(37,1).Value  = 1
$ws.Cells.Item(37,2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(37,3).Value  = "Arica y Parinacota"
$ws.Cells.Item(37,4).Value  = 44414
$ws.Cells.Item(37,5).Value  = 15
$ws.Cells.Item(37,6).Value  = 100114014
$ws.Cells.Item(37,7).Value  = "Betarraga"
$ws.Cells.Item(37,8).Value  = "Sin especificar"
$ws.Cells.Item(37,9).Value  = "Segunda"
$ws.Cells.Item(37,10).Value = 900
$ws.Cells.Item(37,11).Value = 500
$ws.Cells.Item(37,12).Value = 550
$ws.Cells.Item(37,13).Value = 525
$ws.Cells.Item(37,14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(37,15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(37,16).Value = 105
$ws.Cells.Item(37,17).Value = 5
$ws.Cells.Item(37,18).Value = "Hortaliza"
